# Duplicate the "generate-weights-reallife" sheet, placing the copy right
# before the original, then rename the copy to "generate-weights-reallife_2"
# and update its demo data (columns C/D) with new values, dropping the last
# row so the sheet now spans A1:D6 instead of A1:D7.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("generate-weights-reallife")
$source.Copy($source)

# Excel makes the freshly-created copy the active sheet.
$newSheet = $wb.ActiveSheet
$newSheet.Name = "generate-weights-reallife_2"

$newSheet.Range("C3").Value = 0.3
$newSheet.Range("D3").Value = 0.3
$newSheet.Range("C4").Value = 0.5
$newSheet.Range("D4").Value = 0.5
$newSheet.Range("C5").Value = 0.7
$newSheet.Range("D5").Value = 0.7
$newSheet.Range("C6").Value = 1
$newSheet.Range("D6").Value = 1

$newSheet.Range("B7:D7").Delete() | Out-Null

$newSheet.Range("D3").Select() | Out-Null
